# Weekly price-list update: a new week's record is inserted at row 105
# (new date 44512) and every existing record for rows 105-200 shifts
# down by one row, with the former last row (200) landing in the new
# row 201. The sheet's dimension grows from A1:R200 to A1:R201.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 200..105 down into 201..106, working bottom-up so we never
# overwrite a source row before it has been copied.
for ($r = 200; $r -ge 105; $r--) {
    $srcRange = $ws.Range("A" + $r + ":R" + $r)
    $dstRange = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $dstRange.Value2 = $srcRange.Value2
}

# The brand-new row 201 has no pre-existing formatting, so bring over the
# date format used by column D (the only column with non-default styling).
$ws.Range("D201").NumberFormat = $ws.Range("D200").NumberFormat

# Finally, write the new week's date into row 105 (the rest of row 105's
# data - volume/price/unit/origin/etc. - is unchanged).
$ws.Range("D105").Value2 = 44512
